$wb = $excel.ActiveWorkbook

# --- Overview sheet: update combined per-locale status text for the two
#     files that moved from "Ready for handoff" to "In Translation"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: Status column (C) for the same two files
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet: Status column (C) for the same two files
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
